{"js": "// Replace the two-digit multiplication equations throughout the document\n// body (including inside the table cells) with their updated values, per\n// the authoring diff. Each old equation string is unique in the document,\n// so a plain-text search/replace per pair is unambiguous.\nconst replacements = [\n  [\"12\u00d760=720\", \"57\u00d714=798\"],\n  [\"98\u00d714=1372\", \"93\u00d757=5301\"],\n  [\"13\u00d740=520\", \"62\u00d749=3038\"],\n  [\"37\u00d784=3108\", \"89\u00d749=4361\"],\n  [\"16\u00d795=1520\", \"33\u00d777=2541\"],\n  [\"35\u00d741=1435\", \"34\u00d798=3332\"],\n  [\"71\u00d799=7029\", \"39\u00d750=1950\"],\n  [\"90\u00d746=4140\", \"18\u00d788=1584\"],\n  [\"83\u00d738=3154\", \"20\u00d740=800\"],\n  [\"44\u00d761=2684\", \"55\u00d758=3190\"],\n  [\"34\u00d716=544\", \"67\u00d775=5025\"],\n  [\"80\u00d783=6640\", \"93\u00d740=3720\"],\n  [\"59\u00d756=3304\", \"73\u00d754=3942\"],\n  [\"66\u00d773=4818\", \"77\u00d750=3850\"],\n  [\"64\u00d760=3840\", \"97\u00d743=4171\"],\n  [\"41\u00d780=3280\", \"52\u00d790=4680\"],\n  [\"95\u00d757=5415\", \"49\u00d731=1519\"],\n  [\"90\u00d730=2700\", \"44\u00d735=1540\"],\n  [\"62\u00d756=3472\", \"38\u00d764=2432\"],\n  [\"18\u00d727=486\", \"64\u00d743=2752\"],\n  [\"59\u00d728=1652\", \"30\u00d789=2670\"],\n  [\"13\u00d746=598\", \"95\u00d761=5795\"],\n  [\"65\u00d711=715\", \"16\u00d733=528\"],\n  [\"90\u00d785=7650\", \"67\u00d727=1809\"],\n  [\"57\u00d796=5472\", \"21\u00d720=420\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication equations throughout the document\n# (including inside the table cells) with their updated values, per the\n# authoring diff. Each old equation string is unique in the document, so a\n# plain-text Find/Replace per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"12\u00d760=720\",  \"57\u00d714=798\"),\n    @(\"98\u00d714=1372\", \"93\u00d757=5301\"),\n    @(\"13\u00d740=520\",  \"62\u00d749=3038\"),\n    @(\"37\u00d784=3108\", \"89\u00d749=4361\"),\n    @(\"16\u00d795=1520\", \"33\u00d777=2541\"),\n    @(\"35\u00d741=1435\", \"34\u00d798=3332\"),\n    @(\"71\u00d799=7029\", \"39\u00d750=1950\"),\n    @(\"90\u00d746=4140\", \"18\u00d788=1584\"),\n    @(\"83\u00d738=3154\", \"20\u00d740=800\"),\n    @(\"44\u00d761=2684\", \"55\u00d758=3190\"),\n    @(\"34\u00d716=544\",  \"67\u00d775=5025\"),\n    @(\"80\u00d783=6640\", \"93\u00d740=3720\"),\n    @(\"59\u00d756=3304\", \"73\u00d754=3942\"),\n    @(\"66\u00d773=4818\", \"77\u00d750=3850\"),\n    @(\"64\u00d760=3840\", \"97\u00d743=4171\"),\n    @(\"41\u00d780=3280\", \"52\u00d790=4680\"),\n    @(\"95\u00d757=5415\", \"49\u00d731=1519\"),\n    @(\"90\u00d730=2700\", \"44\u00d735=1540\"),\n    @(\"62\u00d756=3472\", \"38\u00d764=2432\"),\n    @(\"18\u00d727=486\",  \"64\u00d743=2752\"),\n    @(\"59\u00d728=1652\", \"30\u00d789=2670\"),\n    @(\"13\u00d746=598\",  \"95\u00d761=5795\"),\n    @(\"65\u00d711=715\",  \"16\u00d733=528\"),\n    @(\"90\u00d785=7650\", \"67\u00d727=1809\"),\n    @(\"57\u00d796=5472\", \"21\u00d720=420\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
